$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of the goods list (column A), starting at row 15 through row 56.
# Row 1-14 (Counts header plus the top 13 goods) are unchanged.
# Only the text values change; the numeric counts in column B stay with their row.
$values = @(
    "железо",
    "Крымскую соль",
    "колеса",
    "полотно",
    "сено",
    "парча",
    "говядина",
    "табак",
    "позумент",
    "выбойка",
    "сахар",
    "шелк",
    "чулок",
    "лес",
    "лыко",
    "китайка",
    "сапог",
    "ладан",
    "сани",
    "коса",
    "ром",
    "горшок",
    "гвоздь",
    "овца",
    "веревка",
    "обод",
    "замок",
    "рогожа",
    "конь",
    "платок",
    "гумми",
    "дуга",
    "брусья",
    "сосуд",
    "бечева",
    "сковорода",
    "покроми",
    "хомут",
    "нитка",
    "роза",
    "скотский кожа",
    "котел"
)

$startRow = 15
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
